# Rename the existing "总计" sheet to "2022-Q1" and give it fresh
# fund-holding data (same layout as the 2020-Q4 / 2021-Q1 / 2021-Q2 sheets).
$wb = $excel.ActiveWorkbook

$q1sheet = $wb.Worksheets.Item("总计")

# Drop the old summary rows (rows 3-4) that don't belong to the new layout;
# row 2 will be overwritten below with the new fund data.
$q1sheet.Range("A3:D4").EntireRow.Delete()

$q1sheet.Name = "2022-Q1"

# Header row
$q1sheet.Range("B1").Value = "基金代码"
$q1sheet.Range("C1").Value = "基金名称"
$q1sheet.Range("D1").Value = "基金规模"
$q1sheet.Range("E1").Value = "股票总仓位"
$q1sheet.Range("F1").Value = "仓位占比"
$q1sheet.Range("G1").Value = "持有市值(亿元)"
$q1sheet.Range("H1").Value = "仓位排名"

# Copy the header style (B1:D1 already carry it) onto the newly added headers
$q1sheet.Range("B1").Copy() | Out-Null
$q1sheet.Range("E1:H1").PasteSpecial(-4122) | Out-Null

# Data row
$q1sheet.Range("A2").Value = 0
$q1sheet.Range("B2").Value = "'006105"
$q1sheet.Range("C2").Value = "泰达宏利印度机会股票（QDII）"
$q1sheet.Range("D2").Value = "'0.60"
$q1sheet.Range("E2").Value = "'87.31"
$q1sheet.Range("F2").Value = "'9.83"
$q1sheet.Range("G2").Value = "'0.0590"
$q1sheet.Range("H2").Value = 1

# Add a brand new "总计" sheet right after "2022-Q1", restoring the aggregate
# view. Duplicate "2022-Q1" first so the new sheet inherits the same sheetPr /
# pageMargins / sheetView layout used across the workbook, then strip it back
# down to the 4-column summary shape.
$q1sheet.Copy($null, $q1sheet)
$totalSheet = $wb.Worksheets.Item($q1sheet.Index + 1)
$totalSheet.Name = "总计"
$totalSheet.Cells.ClearContents()
$totalSheet.Range("E1:H2").Clear()

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.06

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q2"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.04

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.04

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2020-Q4"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.03

# The A2 index cell already carries the header-row style (s=2); replicate it
# down through A3:A5 to match the other "index column" cells.
$totalSheet.Range("A2").Copy() | Out-Null
$totalSheet.Range("A3:A5").PasteSpecial(-4122) | Out-Null

# Restore the originally active sheet/selection so we don't leave an
# unrelated view-state change behind.
$wb.Worksheets.Item("2020-Q4").Activate()
